$p = $ppt.ActivePresentation

# Slide 4, Shape 2 (Content Placeholder 2)
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$tr.Characters(110, 19).Font.Size = 12
$tr.Characters(130, 4).Font.Size = 12
$tr.Characters(136, 6).Font.Size = 12

# Slide 4, Shape 3 (Content Placeholder 3)
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange
$tr.Characters(106, 8).Font.Size = 12
$tr.Characters(148, 1).Font.Size = 12

# Slide 6, Shape 2 (Text Placeholder 3)
$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$tr.Characters(84, 6).Font.Size = 12

# Slide 7, Shape 2 (Text Placeholder 3)
$s = $p.Slides.Item(7)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$tr.Characters(13, 4).Font.Size = 12
$tr.Characters(64, 6).Font.Size = 12
$tr.Characters(81, 4).Font.Size = 12

# Slide 10, Shape 2 (Content Placeholder 2)
$s = $p.Slides.Item(10)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$tr.Characters(29, 4).Font.Size = 12
$tr.Characters(35, 6).Font.Size = 12
$tr.Characters(46, 4).Font.Size = 12
$tr.Characters(216, 4).Font.Size = 12
$tr.Characters(222, 4).Font.Size = 12
$tr.Characters(228, 4).Font.Size = 12
$tr.Characters(238, 3).Font.Size = 12

# Slide 10, Shape 3 (Content Placeholder 3)
$s = $p.Slides.Item(10)
$shp = $s.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, 3).Font.Size = 12
$tr.Characters(5, 5).Font.Size = 12
$tr.Characters(10, 1).Font.Size = 12
$tr.Characters(11, 1).Font.Size = 12
$tr.Characters(12, 29).Font.Size = 12
$tr.Characters(42, 6).Font.Size = 12
$tr.Characters(48, 1).Font.Size = 12
$tr.Characters(49, 1).Font.Size = 12
$tr.Characters(51, 4).Font.Size = 12
$tr.Characters(55, 4).Font.Size = 12
$tr.Characters(59, 1).Font.Size = 12
$tr.Characters(61, 8).Font.Size = 12
$tr.Characters(69, 13).Font.Size = 12
$tr.Characters(82, 1).Font.Size = 12
$tr.Characters(83, 24).Font.Size = 12
$tr.Characters(108, 4).Font.Size = 12
$tr.Characters(112, 8).Font.Size = 12
$tr.Characters(120, 1).Font.Size = 12
$tr.Characters(122, 8).Font.Size = 12
$tr.Characters(130, 11).Font.Size = 12
$tr.Characters(141, 1).Font.Size = 12
$tr.Characters(142, 1).Font.Size = 12
$tr.Characters(143, 5).Font.Size = 12
$tr.Characters(149, 8).Font.Size = 12
$tr.Characters(157, 5).Font.Size = 12
$tr.Characters(162, 1).Font.Size = 12
$tr.Characters(163, 4).Font.Size = 12
$tr.Characters(169, 6).Font.Size = 12
$tr.Characters(175, 1).Font.Size = 12
$tr.Characters(176, 15).Font.Size = 12
$tr.Characters(192, 4).Font.Size = 12
$tr.Characters(196, 1).Font.Size = 12
$tr.Characters(197, 11).Font.Size = 12
$tr.Characters(209, 3).Font.Size = 12

# Slide 11, Shape 2 (Content Placeholder 2)
$s = $p.Slides.Item(11)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, 28).Font.Size = 12
$tr.Characters(31, 34).Font.Size = 12
$tr.Characters(65, 1).Font.Size = 12
$tr.Characters(66, 8).Font.Size = 12
$tr.Characters(74, 1).Font.Size = 12
$tr.Characters(75, 1).Font.Size = 12
$tr.Characters(78, 41).Font.Size = 12
$tr.Characters(120, 44).Font.Size = 12
$tr.Characters(166, 45).Font.Size = 12
$tr.Characters(212, 14).Font.Size = 12
$tr.Characters(226, 9).Font.Size = 12
$tr.Characters(235, 5).Font.Size = 12
$tr.Characters(240, 10).Font.Size = 12
$tr.Characters(250, 4).Font.Size = 12
$tr.Characters(255, 22).Font.Size = 12
$tr.Characters(279, 2).Font.Size = 12
$tr.Characters(282, 8).Font.Size = 12
$tr.Characters(291, 2).Font.Size = 12

# Slide 11, Shape 3 (Content Placeholder 3)
$s = $p.Slides.Item(11)
$shp = $s.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange
$tr.Characters(217, 7).Font.Size = 12
$tr.Characters(229, 8).Font.Size = 12

# Slide 12, Shape 2 (Content Placeholder 2)
$s = $p.Slides.Item(12)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$tr.Characters(121, 6).Font.Size = 12
$tr.Characters(127, 7).Font.Size = 12
$tr.Characters(134, 2).Font.Size = 12
$tr.Characters(136, 3).Font.Size = 12
$tr.Characters(140, 6).Font.Size = 12
$tr.Characters(146, 19).Font.Size = 12
$tr.Characters(165, 2).Font.Size = 12
$tr.Characters(167, 4).Font.Size = 12
$tr.Characters(173, 2).Font.Size = 12
$tr.Characters(175, 1).Font.Size = 12
$tr.Characters(176, 11).Font.Size = 12
$tr.Characters(187, 1).Font.Size = 12
$tr.Characters(188, 2).Font.Size = 12
$tr.Characters(190, 1).Font.Size = 12
$tr.Characters(191, 2).Font.Size = 12
$tr.Characters(193, 4).Font.Size = 12
$tr.Characters(197, 1).Font.Size = 12
$tr.Characters(199, 6).Font.Size = 12
$tr.Characters(205, 1).Font.Size = 12
$tr.Characters(206, 1).Font.Size = 12
$tr.Characters(207, 1).Font.Size = 12
$tr.Characters(208, 1).Font.Size = 12
$tr.Characters(209, 1).Font.Size = 12
$tr.Characters(210, 7).Font.Size = 12
$tr.Characters(217, 1).Font.Size = 12
$tr.Characters(218, 2).Font.Size = 12
$tr.Characters(221, 8).Font.Size = 12
$tr.Characters(229, 1).Font.Size = 12
$tr.Characters(230, 24).Font.Size = 12
$tr.Characters(254, 2).Font.Size = 12
$tr.Characters(257, 17).Font.Size = 12
$tr.Characters(274, 12).Font.Size = 12
$tr.Characters(286, 2).Font.Size = 12
$tr.Characters(288, 7).Font.Size = 12
$tr.Characters(295, 2).Font.Size = 12
$tr.Characters(298, 17).Font.Size = 12
$tr.Characters(316, 15).Font.Size = 12
$tr.Characters(331, 3).Font.Size = 12
$tr.Characters(334, 2).Font.Size = 12
$tr.Characters(336, 1).Font.Size = 12
$tr.Characters(337, 2).Font.Size = 12
$tr.Characters(339, 3).Font.Size = 12
$tr.Characters(342, 2).Font.Size = 12
$tr.Characters(344, 1).Font.Size = 12
$tr.Characters(345, 2).Font.Size = 12
$tr.Characters(348, 8).Font.Size = 12
$tr.Characters(356, 4).Font.Size = 12
$tr.Characters(360, 1).Font.Size = 12
$tr.Characters(362, 10).Font.Size = 12

# Slide 14, Shape 2 (Content Placeholder 2)
$s = $p.Slides.Item(14)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, 38).Font.Size = 12
$tr.Characters(40, 5).Font.Size = 12
$tr.Characters(45, 3).Font.Size = 12
$tr.Characters(48, 14).Font.Size = 12
$tr.Characters(62, 4).Font.Size = 12
$tr.Characters(66, 14).Font.Size = 12
$tr.Characters(81, 5).Font.Size = 12
$tr.Characters(86, 3).Font.Size = 12
$tr.Characters(89, 1).Font.Size = 12
$tr.Characters(92, 46).Font.Size = 12
$tr.Characters(140, 17).Font.Size = 12
$tr.Characters(159, 2).Font.Size = 12
$tr.Characters(161, 36).Font.Size = 12
$tr.Characters(198, 6).Font.Size = 12
$tr.Characters(204, 20).Font.Size = 12
$tr.Characters(225, 6).Font.Size = 12
$tr.Characters(231, 6).Font.Size = 12
$tr.Characters(238, 6).Font.Size = 12
$tr.Characters(244, 8).Font.Size = 12
$tr.Characters(253, 6).Font.Size = 12
$tr.Characters(259, 6).Font.Size = 12
$tr.Characters(266, 6).Font.Size = 12
$tr.Characters(272, 4).Font.Size = 12

# Slide 14, Shape 3 (Content Placeholder 3)
$s = $p.Slides.Item(14)
$shp = $s.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange
$tr.Characters(44, 1).Font.Size = 12
$tr.Characters(59, 2).Font.Size = 12
$tr.Characters(80, 1).Font.Size = 12

# Slide 15, Shape 2 (Content Placeholder 2)
$s = $p.Slides.Item(15)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$tr.Characters(291, 6).Font.Size = 12
$tr.Characters(297, 1).Font.Size = 12
$tr.Characters(298, 13).Font.Size = 12
$tr.Characters(311, 16).Font.Size = 12
$tr.Characters(327, 13).Font.Size = 12
$tr.Characters(340, 1).Font.Size = 12
$tr.Characters(341, 6).Font.Size = 12
$tr.Characters(347, 1).Font.Size = 12
$tr.Characters(349, 7).Font.Size = 12
$tr.Characters(356, 1).Font.Size = 12
$tr.Characters(357, 13).Font.Size = 12
$tr.Characters(370, 26).Font.Size = 12
$tr.Characters(396, 13).Font.Size = 12
$tr.Characters(409, 1).Font.Size = 12
$tr.Characters(410, 6).Font.Size = 12
$tr.Characters(416, 1).Font.Size = 12
$tr.Characters(419, 6).Font.Size = 12
$tr.Characters(425, 1).Font.Size = 12
$tr.Characters(426, 2).Font.Size = 12
$tr.Characters(428, 4).Font.Size = 12
$tr.Characters(432, 2).Font.Size = 12
$tr.Characters(434, 4).Font.Size = 12
$tr.Characters(438, 2).Font.Size = 12
$tr.Characters(440, 4).Font.Size = 12
$tr.Characters(444, 2).Font.Size = 12
$tr.Characters(446, 4).Font.Size = 12
$tr.Characters(450, 2).Font.Size = 12
$tr.Characters(452, 4).Font.Size = 12
$tr.Characters(456, 1).Font.Size = 12
$tr.Characters(458, 8).Font.Size = 12
$tr.Characters(466, 1).Font.Size = 12
$tr.Characters(467, 2).Font.Size = 12
$tr.Characters(469, 9).Font.Size = 12
$tr.Characters(478, 2).Font.Size = 12
$tr.Characters(480, 7).Font.Size = 12
$tr.Characters(487, 2).Font.Size = 12
$tr.Characters(489, 8).Font.Size = 12
$tr.Characters(497, 2).Font.Size = 12
$tr.Characters(499, 7).Font.Size = 12
$tr.Characters(506, 2).Font.Size = 12
$tr.Characters(508, 7).Font.Size = 12
$tr.Characters(515, 2).Font.Size = 12
$tr.Characters(517, 9).Font.Size = 12
$tr.Characters(526, 1).Font.Size = 12
$tr.Characters(529, 20).Font.Size = 12
$tr.Characters(550, 21).Font.Size = 12
$tr.Characters(571, 37).Font.Size = 12
$tr.Characters(608, 1).Font.Size = 12
$tr.Characters(609, 4).Font.Size = 12
$tr.Characters(613, 1).Font.Size = 12
$tr.Characters(614, 1).Font.Size = 12
$tr.Characters(615, 2).Font.Size = 12
$tr.Characters(618, 21).Font.Size = 12
$tr.Characters(639, 24).Font.Size = 12
$tr.Characters(663, 2).Font.Size = 12
$tr.Characters(666, 21).Font.Size = 12
$tr.Characters(687, 23).Font.Size = 12
$tr.Characters(710, 2).Font.Size = 12
$tr.Characters(713, 47).Font.Size = 12
$tr.Characters(760, 2).Font.Size = 12
$tr.Characters(762, 7).Font.Size = 12
$tr.Characters(770, 21).Font.Size = 12
$tr.Characters(791, 9).Font.Size = 12
$tr.Characters(800, 2).Font.Size = 12
$tr.Characters(803, 21).Font.Size = 12
$tr.Characters(824, 23).Font.Size = 12
$tr.Characters(847, 2).Font.Size = 12
$tr.Characters(850, 49).Font.Size = 12
$tr.Characters(899, 10).Font.Size = 12
$tr.Characters(909, 1).Font.Size = 12
$tr.Characters(910, 2).Font.Size = 12
$tr.Characters(912, 13).Font.Size = 12
$tr.Characters(926, 21).Font.Size = 12
$tr.Characters(947, 9).Font.Size = 12
$tr.Characters(956, 2).Font.Size = 12
$tr.Characters(959, 21).Font.Size = 12
$tr.Characters(980, 9).Font.Size = 12
$tr.Characters(989, 2).Font.Size = 12

# Slide 23, Shape 2 (Content Placeholder 2)
$s = $p.Slides.Item(23)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, 3).Font.Size = 12
$tr.Characters(5, 5).Font.Size = 12
$tr.Characters(10, 1).Font.Size = 12
$tr.Characters(11, 1).Font.Size = 12
$tr.Characters(12, 29).Font.Size = 12
$tr.Characters(42, 6).Font.Size = 12
$tr.Characters(48, 1).Font.Size = 12
$tr.Characters(49, 5).Font.Size = 12
$tr.Characters(55, 15).Font.Size = 12
$tr.Characters(71, 28).Font.Size = 12
$tr.Characters(100, 55).Font.Size = 12
$tr.Characters(156, 33).Font.Size = 12
$tr.Characters(190, 59).Font.Size = 12
$tr.Characters(250, 43).Font.Size = 12
$tr.Characters(294, 45).Font.Size = 12
$tr.Characters(341, 6).Font.Size = 12
$tr.Characters(347, 1).Font.Size = 12
$tr.Characters(348, 15).Font.Size = 12
$tr.Characters(364, 4).Font.Size = 12
$tr.Characters(368, 1).Font.Size = 12
$tr.Characters(369, 11).Font.Size = 12
$tr.Characters(381, 3).Font.Size = 12
